$p = $ppt.ActivePresentation

function Add-Paragraph {
    param($TextFrame, [string]$Text)

    $tr = $TextFrame.TextRange
    if ($tr.Length -eq 0) {
        # First paragraph in an (otherwise empty) text frame.
        $tr.Text = $Text
        $target = $TextFrame.TextRange
    } else {
        $TextFrame.TextRange.InsertAfter("`r" + $Text) | Out-Null
        $full = $TextFrame.TextRange
        $target = $full.Characters($full.Length - $Text.Length + 1, $Text.Length)
    }
    $target.LanguageID = "en-US"
}

function Set-BodyText {
    param($TextFrame, [string[]]$Texts, [int[]]$Levels)

    for ($i = 0; $i -lt $Texts.Length; $i++) {
        Add-Paragraph $TextFrame $Texts[$i]
    }

    # Second pass: only touch IndentLevel for paragraphs that need a
    # non-default (>0) outline level, so level-0 paragraphs keep no
    # <a:pPr> element at all (matching how PowerPoint itself omits it).
    $pos = 1
    for ($i = 0; $i -lt $Texts.Length; $i++) {
        $len = $Texts[$i].Length
        if ($Levels[$i] -gt 0) {
            $full = $TextFrame.TextRange
            $target = $full.Characters($pos, $len)
            $target.IndentLevel = $Levels[$i] + 1
        }
        $pos = $pos + $len + 1
    }
}

# --- Slide 8: "Subsystems - LED Controller" ---------------------------
$s8 = $p.Slides.Item(8)
$tf8 = $s8.Shapes.Item(2).TextFrame

$texts8 = @(
    "Objective: Send data from the state composer to the LED registers",
    "  Create LED color object",
    "Set RX and TX pin modes",
    "If Serial buffer is not empty",
    "Read LED index",
    "Read red value",
    "Read green value",
    "Read blue value",
    "Set RGB values for given LED index"
)
$levels8 = @(0, 1, 1, 1, 2, 2, 2, 2, 2)

Set-BodyText $tf8 $texts8 $levels8

# --- Slide 9: "Struggles" ----------------------------------------------
$s9 = $p.Slides.Item(9)
$tf9 = $s9.Shapes.Item(2).TextFrame

$texts9 = @(
    "Communication between the Arduino Nano and the Raspberry pi",
    "Needed to know type of input to the Nano",
    "Serial read efficiency on Arduino Nano"
)
$levels9 = @(0, 1, 0)

Set-BodyText $tf9 $texts9 $levels9
